$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column), shifting
# "Late" / "Outstanding" headers (and all data below them) one column
# to the right.
$ws.Range("N1").EntireColumn.Insert() | Out-Null

# Give the newly inserted column the same stored width as column M
# (11 "characters"), instead of leaving it at the default width.
$ws.Range("N1").ColumnWidth = $ws.Range("M1").ColumnWidth

# Make "Repayment schedule" the active sheet (it was "Transactions"
# before), and leave the selection on R7 as in the saved file.
$ws.Activate() | Out-Null
$ws.Range("R7").Select() | Out-Null
